$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 260 (shifts existing rows 260:281 down to 261:282)
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with the new weekly price record
$ws.Cells.Item(260, 1).Value = 4
$ws.Cells.Item(260, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(260, 3).Value = "Los Lagos"
$ws.Cells.Item(260, 4).Value = 44578
$ws.Cells.Item(260, 5).Value = 10
$ws.Cells.Item(260, 6).Value = 100114013
$ws.Cells.Item(260, 7).Value = "Zanahoria"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 150
$ws.Cells.Item(260, 11).Value = 12000
$ws.Cells.Item(260, 12).Value = 12000
$ws.Cells.Item(260, 13).Value = 12000
$ws.Cells.Item(260, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(260, 15).Value = "Región de Ñuble"
$ws.Cells.Item(260, 16).Value = 600
$ws.Cells.Item(260, 17).Value = 20
$ws.Cells.Item(260, 18).Value = "Hortaliza"
